# Commit: add a new "ODI Batting Extra" worksheet (4th tab) with a header
# row (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL,
# MAN_OF_MATCH) and one data row for match 4485.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands in the 4th
# (final) tab position, matching "Player Info", "ODI Batting", "ODI Bowling",
# "ODI Batting Extra".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Reuse the existing bold/bordered/centered header style already used by the
# other sheets' header rows (copy formats only from an existing header cell
# so no new style entries are created).
$srcHeaderCell = $wb.Worksheets.Item("Player Info").Range("A1")
$srcHeaderCell.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data row. MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are stored as
# literal text (not re-interpreted as numbers/percentages), BATTING_POSITION
# is numeric, MAN_OF_MATCH is plain text.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("A2") "4485"
$ws.Range("B2").Value = 7
Set-TextValue $ws.Range("C2") "0"
Set-TextValue $ws.Range("D2") "0"
Set-TextValue $ws.Range("E2") "3.11%"
$ws.Range("F2").Value = "NO"
